# Round 101 update for UK panel sheet.
# Adds two new survey rows (100 -> wave 41 panel E, 101 -> wave 41 panel F)
# to the bottom of the data table, mirroring the existing pattern used by
# the shared "spss_name" formula in column I.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# ---------------------------------------------------------------------
# Row 113: survey round 100, panel E, wave 41
# ---------------------------------------------------------------------
$ws.Cells.Item(113, 1).Value = 8
$ws.Cells.Item(113, 2).Value = 0
$ws.Cells.Item(113, 3).Value = "uk"
$ws.Cells.Item(113, 4).Value = 100
$ws.Cells.Item(113, 5).Value = "E"
$ws.Cells.Item(113, 6).Value = 41
$ws.Cells.Item(113, 7).Value = 44617
$ws.Cells.Item(113, 8).Value = "21-088043_PEW41_Final_ICUO"
$ws.Cells.Item(113, 9).Formula = '=C113&"_"&"sr"&TEXT(D113,"00")&"_"&YEAR(G113)&TEXT(G113,"MM")&TEXT(G113,"DD")&"_p"&E113&"_wv"&TEXT(F113,"00")&""'
$ws.Cells.Item(113, 10).Value = 1
$ws.Cells.Item(113, 11).Value = 1

# Match the date format already used in column G (numFmtId 14, m/d/yyyy)
$ws.Range("G112").Copy()
$ws.Range("G113").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Row 114: survey round 101, panel F, wave 41
# ---------------------------------------------------------------------
$ws.Cells.Item(114, 1).Value = 8
$ws.Cells.Item(114, 2).Value = 0
$ws.Cells.Item(114, 3).Value = "uk"
$ws.Cells.Item(114, 4).Value = 101
$ws.Cells.Item(114, 5).Value = "F"
$ws.Cells.Item(114, 6).Value = 41
$ws.Cells.Item(114, 7).Value = 44624
$ws.Cells.Item(114, 8).Value = "21-088071_PFW41_Final_ICUO"
$ws.Cells.Item(114, 9).Formula = '=C114&"_"&"sr"&TEXT(D114,"00")&"_"&YEAR(G114)&TEXT(G114,"MM")&TEXT(G114,"DD")&"_p"&E114&"_wv"&TEXT(F114,"00")&""'
$ws.Cells.Item(114, 10).Value = 1
$ws.Cells.Item(114, 11).Value = 1

$ws.Range("G112").Copy()
$ws.Range("G114").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Recalculate so the cached formula results stored in the file are correct.
$excel.CalculateFull()

# ---------------------------------------------------------------------
# Update the view so it scrolls/selects the same way as the saved file.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 71
$win.ScrollColumn = 1
$ws.Range("H86").Select()
